$wb = $excel.ActiveWorkbook

# --- Summary sheet: refresh "Report Generated" timestamp ---
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B2").Value = "2026-02-18 17:21"

# --- Detail sheet: refresh scan data (versions, hosts, instance ids, dates) ---
$detail = $wb.Worksheets.Item("Detail")

# A few values (e.g. "2.4.11", "2025-09-06") look like dates/numbers to Excel's
# auto-detection, so force those specific cells to Text format before writing,
# matching the plain-text values already used throughout this column in the sheet.
$detail.Range("G2").Value = "2.1.4"
$detail.Range("I2").Value = "ip-10-100-3-123.ap-southeast-1.compute.internal"
$detail.Range("J2").Value = "i-0f3c4563f1c12a574"
$detail.Range("G3").Value = "2.3.1"
$detail.Range("H3").Value = "2.5.4"
$detail.Range("I3").Value = "ip-172-17-2-49.ap-southeast-1.compute.internal"
$detail.Range("J3").Value = "i-0799e3a4b56be2630"
$detail.Range("G4").Value = "4.0.1"
$detail.Range("H4").Value = "4.0.4"
$detail.Range("I4").Value = "ip-10-100-3-123.ap-southeast-1.compute.internal"
$detail.Range("J4").Value = "i-0f3c4563f1c12a574"
$detail.Range("I5").Value = "ip-172-17-2-49.ap-southeast-1.compute.internal"
$detail.Range("J5").Value = "i-0799e3a4b56be2630"
$detail.Range("G15").Value = "3.1.2"
$detail.Range("I15").Value = "ip-10-100-3-123.ap-southeast-1.compute.internal"
$detail.Range("J15").Value = "i-0f3c4563f1c12a574"
$detail.Range("K15").NumberFormat = "@"
$detail.Range("K15").Value = "2025-09-06"
$detail.Range("G16").Value = "3.0.14"
$detail.Range("I16").Value = "ip-172-17-2-49.ap-southeast-1.compute.internal"
$detail.Range("J16").Value = "i-0799e3a4b56be2630"
$detail.Range("K16").NumberFormat = "@"
$detail.Range("K16").Value = "2025-06-26"
$detail.Range("G48").Value = "4.17.5"
$detail.Range("I48").Value = "ip-172-17-2-49.ap-southeast-1.compute.internal"
$detail.Range("J48").Value = "i-0799e3a4b56be2630"
$detail.Range("K48").NumberFormat = "@"
$detail.Range("K48").Value = "2025-03-18"
$detail.Range("G49").Value = "4.17.4"
$detail.Range("G50").Value = "3.10.1"
$detail.Range("I50").Value = "ip-10-100-3-123.ap-southeast-1.compute.internal"
$detail.Range("J50").Value = "i-0f3c4563f1c12a574"
$detail.Range("K50").NumberFormat = "@"
$detail.Range("K50").Value = "2025-05-31"
$detail.Range("G52").NumberFormat = "@"
$detail.Range("G52").Value = "2.4.11"
$detail.Range("I52").Value = "ip-10-100-3-123.ap-southeast-1.compute.internal"
$detail.Range("J52").Value = "i-0f3c4563f1c12a574"
$detail.Range("K52").NumberFormat = "@"
$detail.Range("K52").Value = "2025-09-06"
$detail.Range("G53").Value = "2.4.8"
$detail.Range("K53").NumberFormat = "@"
$detail.Range("K53").Value = "2025-08-24"
$detail.Range("G54").NumberFormat = "@"
$detail.Range("G54").Value = "2.4.10"
$detail.Range("I54").Value = "ip-172-17-2-49.ap-southeast-1.compute.internal"
$detail.Range("J54").Value = "i-0799e3a4b56be2630"
$detail.Range("I56").Value = "ip-172-17-2-49.ap-southeast-1.compute.internal"
$detail.Range("J56").Value = "i-0799e3a4b56be2630"
$detail.Range("I57").Value = "ip-10-100-3-123.ap-southeast-1.compute.internal"
$detail.Range("J57").Value = "i-0f3c4563f1c12a574"
